# Scenario step in the "Use Case" section originally read:
#   "Return a book that is one day overdue"
# It should become:
#   "Return a book that is two days overdue"
#
# Locate the paragraph by its exact original text (rather than a fixed
# index) so the script is resilient to any other paragraph-count
# differences, then perform two tightly-scoped, whole-word Find &
# Replace operations restricted to that paragraph's own Range so that
# the many other "day"/"days" mentions elsewhere in the document are
# left untouched.

$d = $word.ActiveDocument

$oldText = "Return a book that is one day overdue"
$target = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq ($oldText + [char]13)) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the 'Return a book that is one day overdue' paragraph"
}

# "one" -> "two"
$rng1 = $target.Range
$rng1.Find.ClearFormatting()
$rng1.Find.Execute("one", $true, $true, $false, $false, $false, $true, 1, $false, "two", 2)

# "day" -> "days" (re-fetch the range; the paragraph grew by one character)
$rng2 = $target.Range
$rng2.Find.ClearFormatting()
$rng2.Find.Execute("day", $true, $true, $false, $false, $false, $true, 1, $false, "days", 2)

Write-Output $target.Range.Text
